$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 1_Tin_Hieu_Hom_Nay -> add column H "GTGD_TB_Tỷ" ----
$ws1 = $wb.Worksheets.Item("1_Tin_Hieu_Hom_Nay")

# Copy formatting from header style cell (A1) to new header H1, then set its text
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws1.Range("H1").Value = "GTGD_TB_Tỷ"

$h1Values = @{
    2 = 16.369
    3 = 1.715
    4 = 0.351
    5 = 0.196
    6 = 1.527
    7 = 0.995
    8 = 1.668
    9 = 2.442
    10 = 20.628
    11 = 12.683
    12 = 0.646
    13 = 1.507
    14 = 3.696
    15 = 0.33
    16 = 13.32
    17 = 3.991
    18 = 1.819
    19 = 7.742
    20 = 4.737
    21 = 6.824
    22 = 14.76
    23 = 15.044
    24 = 12.044
    25 = 1.319
    26 = 1.355
    27 = 0.499
    28 = 5.871
    29 = 0.14
    30 = 7.097
    31 = 1.997
    32 = 8.959
    33 = 0.796
    34 = 1.746
    35 = 0.653
    36 = 11.833
    37 = 6.621
    38 = 1.29
    39 = 0.506
    40 = 62.023
    41 = 1.948
    42 = 0.095
    43 = 0.836
    44 = 6.209
    45 = 5.834
    46 = 6.747
    47 = 4.598
    48 = 1.162
    49 = 2.088
    50 = 4.343
    51 = 3.937
    52 = 18.708
}
foreach ($r in $h1Values.Keys) {
    $ws1.Cells.Item($r, 8).Value = $h1Values[$r]
}

# ---- Sheet 2: 2_Xu_Huong_21_Ngay -> add column G "GTGD_TB_Tỷ" ----
$ws2 = $wb.Worksheets.Item("2_Xu_Huong_21_Ngay")

$ws2.Range("A1").Copy() | Out-Null
$ws2.Range("G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws2.Range("G1").Value = "GTGD_TB_Tỷ"

$g1Values = @{
    2 = 13.32
    3 = 12.044
    4 = 3.991
    5 = 20.628
    6 = 12.683
    7 = 18.708
    8 = 6.209
    9 = 5.834
    10 = 3.937
    11 = 4.737
    12 = 14.76
    13 = 4.598
    14 = 4.343
    15 = 8.959
    16 = 16.369
    17 = 6.621
    18 = 15.044
    19 = 1.507
    20 = 7.742
    21 = 6.824
    22 = 1.715
    23 = 2.442
    24 = 62.023
    25 = 0.196
    26 = 1.162
    27 = 0.796
    28 = 1.668
    29 = 7.097
    30 = 5.871
    31 = 0.506
    32 = 2.088
    33 = 1.948
    34 = 1.819
    35 = 0.653
    36 = 1.29
    37 = 0.995
    38 = 0.095
    39 = 1.746
    40 = 3.696
    41 = 1.319
    42 = 0.351
    43 = 11.833
    44 = 1.997
    45 = 0.646
    46 = 0.33
    47 = 0.836
    48 = 1.355
    49 = 1.527
    50 = 0.499
    51 = 0.14
    52 = 6.747
}
foreach ($r in $g1Values.Keys) {
    $ws2.Cells.Item($r, 7).Value = $g1Values[$r]
}

Write-Host "Update complete"